# "Generate Report for Handoff"
#
# Updates the localization-status report to reflect that the content is now
# ready for handoff (was "In Translation"), and refreshes the associated
# timestamps. Also widens the per-language "status/handoff" columns so the
# longer "Ready for handoff" label isn't truncated.

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item(1)   # "Overview"
$ws_zhcn     = $wb.Worksheets.Item(2)   # "zh-cn"
$ws_dede     = $wb.Worksheets.Item(3)   # "de-de"

# --- Status text: "In Translation" -> "Ready for handoff" -----------------
$ws_overview.Range("E2").Value = "Ready for handoff"
$ws_overview.Range("F2").Value = "Ready for handoff"
$ws_zhcn.Range("C2").Value = "Ready for handoff"
$ws_dede.Range("C2").Value = "Ready for handoff"

# --- Timestamps -------------------------------------------------------------
# Latest HO Xliff Generate Date (Overview G2, de-de H2)
$ws_overview.Range("G2").Value = "2016-08-31 06:44:10"
$ws_dede.Range("H2").Value = "2016-08-31 06:44:10"

# Latest Handoff Datetime (zh-cn H2)
$ws_zhcn.Range("H2").Value = "2016-08-31 06:43:59"

# --- Column widths: widen the status/handoff columns ------------------------
$ws_overview.Columns.Item(5).ColumnWidth = 16.29   # E
$ws_overview.Columns.Item(6).ColumnWidth = 16.29   # F
$ws_zhcn.Columns.Item(3).ColumnWidth = 16.29        # C
$ws_dede.Columns.Item(3).ColumnWidth = 16.29        # C
